# Update grades from emails up to March 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Arevalo, Andres ---
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = 24

# --- Row 4: Baesu, Benjamin ---
$ws.Range("D4").Value = 41
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 24

# --- Row 12: Gibbs, Paul L. ---
$ws.Range("F12").Value = 24

# --- Row 15: Gutierrez, Osvaldo ---
$ws.Range("D15").Value = 35

# --- Row 29: Quijano, Jesse A. ---
$ws.Range("D29").Value = 36
$ws.Range("E29").Value = 17

# Highlight the "points possible" row (row 2) with the accent fill used to
# flag it as the reference/max-points row.
$ws.Range("C2:F2").Interior.ThemeColor = 4
$ws.Range("T2").Interior.ThemeColor = 4

# Restore the active selection to the cell last being edited.
$ws.Range("F13").Select() | Out-Null
